$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Price" column (D) -----------------------------------------------------
# These are plain-text strings in the source data (several use "." as both a
# thousands AND a decimal separator, e.g. "60.325.19", and some have trailing
# zeros that matter, e.g. "7.40"), so a naive .Value assignment would let Excel
# auto-coerce them into numbers/dates and corrupt the text (dropping zeros,
# re-formatting, etc.). Temporarily force Text format for the write, then
# restore the cell's original style so no visible formatting changes.
$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.325.19"
$ws.Range("D2").Style = $style
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.630.36"
$ws.Range("D3").Style = $style
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.01"
$ws.Range("D5").Style = $style
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.26"
$ws.Range("D6").Style = $style
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.608"
$ws.Range("D8").Style = $style
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.654.39"
$ws.Range("D9").Style = $style
$style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.85"
$ws.Range("D10").Style = $style
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.148"
$ws.Range("D12").Style = $style
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.094.87"
$ws.Range("D14").Style = $style
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.297.08"
$ws.Range("D15").Style = $style
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.02"
$ws.Range("D16").Style = $style
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").Style = $style
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.642.60"
$ws.Range("D18").Style = $style
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("D19").Style = $style
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "342.75"
$ws.Range("D20").Style = $style
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.43"
$ws.Range("D21").Style = $style
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("D22").Style = $style
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.91"
$ws.Range("D24").Style = $style
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.439"
$ws.Range("D25").Style = $style
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.40"
$ws.Range("D28").Style = $style
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0806"
$ws.Range("D29").Style = $style
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.15"
$ws.Range("D32").Style = $style
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.51"
$ws.Range("D33").Style = $style
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.18"
$ws.Range("D34").Style = $style
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.11"
$ws.Range("D35").Style = $style
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.889"
$ws.Range("D38").Style = $style
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.66"
$ws.Range("D39").Style = $style
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.52"
$ws.Range("D40").Style = $style
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "300.64"
$ws.Range("D41").Style = $style
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.64"
$ws.Range("D42").Style = $style
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0985"
$ws.Range("D44").Style = $style
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.605"
$ws.Range("D45").Style = $style
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0545"
$ws.Range("D46").Style = $style
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.48"
$ws.Range("D47").Style = $style
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.82"
$ws.Range("D48").Style = $style

# --- "Volume(1h)" column (E) -------------------------------------------------
# Percentage text values (kept as text, with their original padding spaces).
$ws.Range("E2").Value = "  +3.99%  "
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +6.41%  "
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  +4.56%  "
$ws.Range("E9").Value = "  +3.14%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  +6.34%  "
$ws.Range("E12").Value = "  +7.38%  "
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("E17").Value = "  +5.80%  "
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("E20").Value = "  +2.72%  "
$ws.Range("E21").Value = "  +4.59%  "
$ws.Range("E22").Value = "  +4.56%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.93%  "
$ws.Range("E25").Value = "  +5.34%  "
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E28").Value = "  +5.90%  "
$ws.Range("E29").Value = "  +11.72%  "
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  +4.78%  "
$ws.Range("E32").Value = "  +5.56%  "
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("E35").Value = "  +6.50%  "
$ws.Range("E36").Value = "  +7.88%  "
$ws.Range("E37").Value = "  +6.15%  "
$ws.Range("E38").Value = "  +9.56%  "
$ws.Range("E39").Value = "  +2.11%  "
$ws.Range("E40").Value = "  +7.96%  "
$ws.Range("E41").Value = "  +6.72%  "
$ws.Range("E42").Value = "  +2.24%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("E44").Value = "  +4.53%  "
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("E47").Value = "  +7.20%  "
$ws.Range("E48").Value = "  +16.44%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("E50").Value = "  +4.76%  "
$ws.Range("E51").Value = "  +7.78%  "
